$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (title reflects updated "through" date)
$ws.Name = "Through 2021-10-02"

# Update the label text for the October row
$ws.Range("A11").Value = "October (through 10-02)"

# Update October row (row 11) values
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 6
$ws.Range("D11").Value = 8
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 9
$ws.Range("H11").Value = 17

# Update Total row (row 12) values
$ws.Range("B12").Value = 229
$ws.Range("C12").Value = 435
$ws.Range("D12").Value = 635
$ws.Range("E12").Value = 553
$ws.Range("F12").Value = 423
$ws.Range("G12").Value = 910
$ws.Range("H12").Value = 1265
